$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "-test" suffixed values in row 2 to "-test1"
$ws.Range("G2").Value = "abc, ahostess-test1"
$ws.Range("B2").Value = "ahostess-test1 abc"
$ws.Range("C2").Value = "bcohost-test1 abc"
$ws.Range("D2").Value = "guest1-test1 abc"
$ws.Range("E2").Value = "guest2-test1 abc"
$ws.Range("F2").Value = "guest3-test1 abc"
$ws.Range("H2").Value = "abc, bcohost-test1"
$ws.Range("J2").Value = "abc, guest1-test1"
$ws.Range("K2").Value = "abc, guest2-test1"

# Update the view: drop the frozen/scrolled top-left cell and move the
# active selection to D8
$ws.Range("D8").Select()
